# Apply the "Hong Kong EPS v2.0.0" update to
# InputData/trans/EoDfVUwFC/Elast of Demand for Veh Use wrt Fuel Cost.xlsx
#
# Summary of changes:
#  - "EoDfVUwFC" sheet: header B1 becomes "Elasticity (dimensionless)",
#    wraps text, and row 1 grows to fit it.
#  - "About" sheet gains a new "Notes" section (rows 48-50) explaining
#    the rebound-effect elasticity.
#  - Leave behind the same cursor/selection state the authoring session
#    ended on: B1 selected on EoDfVUwFC, A48 selected (and active) on About.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsElast = $wb.Worksheets.Item("EoDfVUwFC")

# ---------------------------------------------------------------------
# About sheet: append a Notes section under the existing content
# ---------------------------------------------------------------------
$wsAbout.Range("A48").Value = "Notes"
$wsAbout.Range("A48").Font.Bold = $true

$wsAbout.Range("A49").Value = 'This variable is also known as the "Fuel Economy Rebound Effect" or "Fuel Cost Rebound Effect." It is the change'
$wsAbout.Range("A50").Value = 'in VMT as a fraction of the change in fuel cost. E.g. for a 1% increase in fuel cost per mile, VMT changes by -0.1%.'

# ---------------------------------------------------------------------
# EoDfVUwFC sheet: rename header and wrap it onto two lines
# ---------------------------------------------------------------------
$wsElast.Range("B1").Value = "Elasticity (dimensionless)"
$wsElast.Range("B1").WrapText = $true
$wsElast.Rows.Item(1).RowHeight = 30

# ---------------------------------------------------------------------
# Restore the on-disk cursor/selection state: EoDfVUwFC was left with
# B1 selected, while About (the visible/active tab) was scrolled down
# and left with A48 selected.
# ---------------------------------------------------------------------
$wsElast.Range("B1").Select()
$wsAbout.Activate()
$wsAbout.Range("A48").Select()
